$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.172.79"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "2.527.49"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'590.56"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'172.98"
$ws.Range("E6").Value = "  +3.74%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "2.526.40"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "2.987.43"
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "67.058.46"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "2.519.03"
$ws.Range("E18").Value = "  -3.04%  "
$ws.Range("E19").Value = "  +4.60%  "
$ws.Range("D20").Value = "'11.38"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "'354.25"
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("E24").Value = "  +5.93%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'69.71"
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").Value = "'9.96"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.653.81"
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("D30").Value = "0.0₃0978"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").Value = "'532.59"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'157.61"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "'18.64"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "'18.45"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +3.11%  "
$ws.Range("D46").Value = "'149.05"
$ws.Range("D47").Value = "'0.557"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("E48").Value = "  -3.09%  "
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("E51").Value = "  -0.40%  "

# Reset style on cells that required a leading apostrophe to stay text,
# so no stray numFmt/quotePrefix style index lingers on the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
